# FN-3460: add two rows previously causing FP errors in facility utilisation cells
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5 - "Crumpet" facility (new row, inserted before the previously-blank row 6)
$ws.Range("A5").Value = "Crumpet GEF"
$ws.Range("B5").Value = 20001371
$ws.Range("C5").Value = "Crumpet exporter"
$ws.Range("D5").Value = "GBP"
$ws.Range("E5").Value = 7000000
$ws.Range("F5").Value = 3938753.8
$ws.Range("G5").Value = 777
$ws.Range("H5").Value = 456
$ws.Range("I5").Value = "GBP"
$ws.Range("J5").Value = "GBP"

# Row 6 - "Scone" facility (previously a blank styled row, now populated)
$ws.Range("A6").Value = "Scone GEF"
$ws.Range("B6").Value = 20001371
$ws.Range("C6").Value = "Scone exporter"
$ws.Range("D6").Value = "GBP"
$ws.Range("E6").Value = 770000
$ws.Range("F6").Value = 761579.37
$ws.Range("G6").Value = 777
$ws.Range("H6").Value = 456.77
$ws.Range("I6").Value = "GBP"
$ws.Range("J6").Value = "GBP"

# Row 5 is a brand-new row with no pre-existing style, so copy the
# formatting of row 4 (A:J) down onto it to match the existing data-row
# style used by rows 2-4. Row 6 was already a styled (but empty) row, so
# its per-cell styles are left as-is and only the values are populated.
$ws.Range("A4:J4").Copy()
$ws.Range("A5:J5").PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = 0

# Update selection to match the new focus of A5:J6 (matches diff's <selection>)
$ws.Range("A5:J6").Select()

$wb.Save()
